# New Submission Synced: 2026-02-08 16:47:09
# The "JSS 3E" sheet is a Google-Forms-style response log (Timestamp,
# Full Name, Admission No, AI Score). A new submission came in, so a new
# row is appended — and, as happens on every sync, the previous last
# row's "Admission No" (which had been left as raw synced text) gets
# normalized to a real number while the row that was *just* synced keeps
# its text representation until the next sync touches it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Normalize the previous row's Admission No ("45") from text to a number.
$ws.Cells.Item(3, 3).Value = 45

# Append the newly synced submission as row 4.
$ws.Cells.Item(4, 1).Value = "2026-02-08 16:47:09"
$ws.Cells.Item(4, 2).Value = "ELISHA BITRUS DAUDA"

# Admission No "7" stays as text for the freshly-synced row (matches the
# sync tool's raw-string behavior for unprocessed rows) — force text so
# Excel doesn't auto-coerce it to a number, then drop back to the default
# style so no stray number-format is left on the cell.
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "7"
$ws.Cells.Item(4, 3).Style = "Normal"

$ws.Cells.Item(4, 4).Value = 10
